# Updated cryptos list on Sat Aug  3 14:00:56 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (col D) / Volume(1h) (col E) updates per row.
# D values are stored as plain text in the sheet (e.g. "542.70", "1.00"),
# so we force text entry with a leading apostrophe and then strip the
# resulting "quote prefix" number-format flag via ClearFormats() so the
# cell's style index is left untouched.
$updates = @(
    @{ Row = 2;  D = "61.841.80"; E = "  -4.79%  " }
    @{ Row = 3;  D = "2.991.96";  E = "  -4.84%  " }
    @{ Row = 4;  D = $null;       E = "  -0.03%  " }
    @{ Row = 5;  D = "542.70";    E = "  -5.53%  " }
    @{ Row = 6;  D = "151.87";    E = "  -8.60%  " }
    @{ Row = 7;  D = $null;       E = "  -0.10%  " }
    @{ Row = 8;  D = "0.570";     E = "  -0.98%  " }
    @{ Row = 9;  D = "3.004.09";  E = "  -5.05%  " }
    @{ Row = 10; D = "0.114";     E = "  -4.26%  " }
    @{ Row = 11; D = $null;       E = "  -7.80%  " }
    @{ Row = 12; D = $null;       E = "  -3.80%  " }
    @{ Row = 13; D = "3.512.10";  E = "  -5.09%  " }
    @{ Row = 14; D = $null;       E = "  -1.18%  " }
    @{ Row = 15; D = "61.863.74"; E = "  -4.78%  " }
    @{ Row = 16; D = "23.97";     E = "  -4.76%  " }
    @{ Row = 17; D = "2.998.41";  E = "  -5.06%  " }
    @{ Row = 18; D = $null;       E = "  -5.82%  " }
    @{ Row = 19; D = $null;       E = "  -1.56%  " }
    @{ Row = 20; D = $null;       E = "  -4.07%  " }
    @{ Row = 21; D = "377.79";    E = "  -8.21%  " }
    @{ Row = 22; D = $null;       E = "  -5.47%  " }
    @{ Row = 23; D = $null;       E = "  +0.35%  " }
    @{ Row = 24; D = "66.17";     E = "  -3.99%  " }
    @{ Row = 25; D = "3.115.03";  E = "  -5.03%  " }
    @{ Row = 26; D = $null;       E = "  -3.43%  " }
    @{ Row = 27; D = $null;       E = "  -3.43%  " }
    @{ Row = 28; D = $null;       E = "  +0.24%  " }
    @{ Row = 29; D = $null;       E = "  -10.23%  " }
    @{ Row = 30; D = "8.25";      E = "  -10.50%  " }
    @{ Row = 31; D = "1.00";      E = "  +0.03%  " }
    @{ Row = 32; D = $null;       E = "  -5.02%  " }
    @{ Row = 33; D = "20.48";     E = "  -4.30%  " }
    @{ Row = 34; D = "160.70";    E = "  -1.94%  " }
    @{ Row = 35; D = $null;       E = "  -4.99%  " }
    @{ Row = 36; D = "4.59";      E = "  -8.18%  " }
    @{ Row = 37; D = $null;       E = "  -6.49%  " }
    @{ Row = 38; D = $null;       E = "  -6.71%  " }
    @{ Row = 39; D = $null;       E = "  -8.19%  " }
    @{ Row = 40; D = "37.49";     E = "  -2.29%  " }
    @{ Row = 41; D = "2.416.69";  E = "  -7.45%  " }
    @{ Row = 42; D = $null;       E = "  -6.32%  " }
    @{ Row = 43; D = "22.03";     E = "  -8.27%  " }
    @{ Row = 44; D = $null;       E = "  -3.45%  " }
    @{ Row = 45; D = $null;       E = "  -5.09%  " }
    @{ Row = 46; D = "5.20";      E = "  -2.40%  " }
    @{ Row = 47; D = "0.996";     E = "  +0.04%  " }
    @{ Row = 48; D = $null;       E = "  -5.00%  " }
    @{ Row = 49; D = "0.0952";    E = "  -2.51%  " }
    @{ Row = 50; D = "266.88";    E = "  -8.13%  " }
    @{ Row = 51; D = "19.62";     E = "  -8.91%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Value = "'" + $u.D
        $cell.ClearFormats()
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
